$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 datetimes (Correspond Handoff Datetime / Correspond Handback DateTime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 14:10:09"
$wsZhCn.Range("G5").Value = "2016-01-25 14:10:50"

# de-de sheet: row 5 datetimes (Correspond Handoff Datetime / Correspond Handback DateTime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 14:10:18"
$wsDeDe.Range("G5").Value = "2016-01-25 14:11:10"
